{"js": "// Remove the empty paragraph, the \"Ver no Jupiter...\" paragraph and the\n// \"\u00a9 2020 ...\" footer paragraph that used to follow the \"LOM3036: ...\"\n// requirement line, mirroring the site-regeneration diff that dropped the\n// trailing \"Ver no Jupiter / copyright\" boilerplate block.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\",\n];\n\n// Find the \"Ver no Jupiter...\" paragraph; the diff also drops the blank\n// paragraph immediately preceding it (right after \"LOM3036: ...\") and the\n// copyright paragraph immediately following it.\nconst items = paragraphs.items;\nlet verIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === targets[0]) {\n    verIdx = i;\n    break;\n  }\n}\n\nif (verIdx !== -1) {\n  const toDelete = [];\n  // Blank paragraph right before \"Ver no Jupiter...\" (only if really blank).\n  if (verIdx - 1 >= 0 && items[verIdx - 1].text === \"\") {\n    toDelete.push(items[verIdx - 1]);\n  }\n  toDelete.push(items[verIdx]);\n  // Copyright paragraph right after, if present.\n  if (verIdx + 1 < items.length && items[verIdx + 1].text === targets[1]) {\n    toDelete.push(items[verIdx + 1]);\n  }\n  for (const p of toDelete) {\n    p.delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" line,\n# the \"\u00a9 2020 ... Creative Commons Attribution\" copyright line, and the blank\n# paragraph between them and the \"LOM3036: ...\" requirement line above them \u2014\n# boilerplate dropped by the site regeneration described in the commit.\n$d = $word.ActiveDocument\n\n$verText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = [char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$count = $d.Paragraphs.Count\n$verIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n  $t = $d.Paragraphs.Item($i).Range.Text\n  $t = $t.TrimEnd([char]13, [char]7)\n  if ($t -eq $verText) {\n    $verIdx = $i\n    break\n  }\n}\n\nif ($verIdx -ge 1) {\n  # Delete highest index first so earlier indices stay valid.\n  if ($verIdx + 1 -le $d.Paragraphs.Count) {\n    $afterText = $d.Paragraphs.Item($verIdx + 1).Range.Text.TrimEnd([char]13, [char]7)\n    if ($afterText -eq $copyrightText) {\n      $d.Paragraphs.Item($verIdx + 1).Range.Delete()\n    }\n  }\n\n  $d.Paragraphs.Item($verIdx).Range.Delete()\n\n  if ($verIdx - 1 -ge 1) {\n    $beforeText = $d.Paragraphs.Item($verIdx - 1).Range.Text.TrimEnd([char]13, [char]7)\n    if ($beforeText -eq \"\") {\n      $d.Paragraphs.Item($verIdx - 1).Range.Delete()\n    }\n  }\n}\n"}
